$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Total Evening" label and SUM formula for evening shift (rows 6:9)
# (entered first so it lands at shared-string index 8)
$ws.Range("E7").Value = "Total Evening"
$ws.Range("F7").Formula = "=SUM(C6:C9)"

# Add "Total morning" label and SUM formula for morning shift (rows 2:5)
# (entered second so it lands at shared-string index 9)
$ws.Range("E3").Value = "Total morning"
$ws.Range("F3").Formula = "=SUM(C2:C5)"

# Set column E width (target stored width ~14.42578125 chars; engine quantizes
# ColumnWidth to 1/6-character steps, so 13.6 is the closest achievable input)
$ws.Columns.Item(5).ColumnWidth = 13.6

# Update selection to F4
$ws.Range("F4").Select()
